$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 93 (new record for date 2021-02-08 = 44235),
# pushing the former rows 93..113 down to 94..114.
$ws.Rows.Item(93).Insert()

# Fix the inserted row's formatting on column A (date column) to match the
# surrounding rows - Insert() otherwise leaves a stray bold/format style.
$ws.Range("A92").Copy()
$ws.Range("A93").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the rolling-window values that shifted because of the newly
# inserted record (rows 90-96 and 112 recompute; 97-111 are unaffected).
$ws.Range("C90").Value = 18
$ws.Range("D90").Value = 274.1395065488882

$ws.Range("C91").Value = 20
$ws.Range("D91").Value = 304.5994517209869

$ws.Range("C92").Value = 19
$ws.Range("D92").Value = 289.3694791349375

$ws.Range("A93").Value = 44235
$ws.Range("B93").Value = 2
$ws.Range("C93").Value = 19
$ws.Range("D93").Value = 289.3694791349375

$ws.Range("C94").Value = 16
$ws.Range("D94").Value = 243.6795613767895

$ws.Range("C95").Value = 13
$ws.Range("D95").Value = 197.9896436186415

$ws.Range("C96").Value = 13
$ws.Range("D96").Value = 197.9896436186415

$ws.Range("C112").Value = 31
$ws.Range("D112").Value = 472.1291501675296

# Append a brand-new last row (115) for date 2021-03-02.
# (C115/D115 stay blank, same as the other "not enough data yet" rows.)
$ws.Range("A115").Value = 44257
$ws.Range("B115").Value = 6

$ws.Range("A114").Copy()
$ws.Range("A115").PasteSpecial(-4122)
$excel.CutCopyMode = $false
